$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '26.616.83'
$ws.Range("E2").Value = '  -0.18%  '
$ws.Range("D3").Value = '1.595.26'
$ws.Range("E3").Value = '  +0.21%  '
$ws.Range("E4").Value = '  +0.00%  '
$c = $ws.Range("D5")
$c.NumberFormat = "@"
$c.Value = '211.03'
$c.Style = "Normal"
$ws.Range("E5").Value = '  -0.01%  '
$ws.Range("E6").Value = '  +1.02%  '
$ws.Range("E7").Value = '  +0.03%  '
$ws.Range("E8").Value = '  -0.09%  '
$ws.Range("E9").Value = '  -1.57%  '
$c = $ws.Range("D10")
$c.NumberFormat = "@"
$c.Value = '19.37'
$c.Style = "Normal"
$ws.Range("E10").Value = '  -1.66%  '
$c = $ws.Range("D11")
$c.NumberFormat = "@"
$c.Value = '0.0836'
$c.Style = "Normal"
$ws.Range("E11").Value = '  +0.20%  '
$ws.Range("D12").Value = '1.818.62'
$ws.Range("E12").Value = '  +0.22%  '
$ws.Range("D13").Value = '1.566.46'
$ws.Range("E13").Value = '  -1.51%  '
$ws.Range("E14").Value = '  -0.23%  '
$ws.Range("E15").Value = '  -1.26%  '
$c = $ws.Range("D16")
$c.NumberFormat = "@"
$c.Value = '64.61'
$c.Style = "Normal"
$ws.Range("E16").Value = '  -0.15%  '
$ws.Range("D17").Value = '26.589.72'
$ws.Range("E17").Value = '  -0.23%  '
$ws.Range("D18").Value = '0.0₃0730'
$ws.Range("E18").Value = '  +0.42%  '
$ws.Range("E19").Value = '  +0.05%  '
$c = $ws.Range("D20")
$c.NumberFormat = "@"
$c.Value = '207.68'
$c.Style = "Normal"
$ws.Range("E20").Value = '  -0.15%  '
$c = $ws.Range("D21")
$c.NumberFormat = "@"
$c.Value = '6.90'
$c.Style = "Normal"
$ws.Range("E21").Value = '  +2.11%  '
$ws.Range("E22").Value = '  -0.14%  '
$ws.Range("E23").Value = '  -3.46%  '
$c = $ws.Range("D24")
$c.NumberFormat = "@"
$c.Value = '8.86'
$c.Style = "Normal"
$ws.Range("E24").Value = '  -0.71%  '
$c = $ws.Range("D25")
$c.NumberFormat = "@"
$c.Value = '145.39'
$c.Style = "Normal"
$ws.Range("E25").Value = '  -1.14%  '
$ws.Range("E26").Value = '  +0.13%  '
$ws.Range("E27").Value = '  -2.02%  '
$ws.Range("E28").Value = '  +0.04%  '
$ws.Range("E29").Value = '  -0.19%  '
$ws.Range("E30").Value = '  -0.28%  '
$c = $ws.Range("D31")
$c.NumberFormat = "@"
$c.Value = '1.15'
$c.Style = "Normal"
$ws.Range("E31").Value = '  +0.02%  '
$ws.Range("E32").Value = '  -0.16%  '
$ws.Range("E33").Value = '  +0.54%  '
$ws.Range("E34").Value = '  +0.16%  '
$ws.Range("D35").Value = '1.284.77'
$ws.Range("E35").Value = '  -2.24%  '
$ws.Range("E36").Value = '  +1.60%  '
$ws.Range("E37").Value = '  -0.51%  '
$ws.Range("E38").Value = '  -0.54%  '
$ws.Range("E39").Value = '  +1.10%  '
$ws.Range("E40").Value = '  +0.07%  '
$ws.Range("E41").Value = '  +0.67%  '
$ws.Range("E42").Value = '  +1.06%  '
$c = $ws.Range("D43")
$c.NumberFormat = "@"
$c.Value = '0.785'
$c.Style = "Normal"
$ws.Range("E43").Value = '  -0.69%  '
$c = $ws.Range("D44")
$c.NumberFormat = "@"
$c.Value = '63.60'
$c.Style = "Normal"
$ws.Range("E44").Value = '  +0.19%  '
$c = $ws.Range("D45")
$c.NumberFormat = "@"
$c.Value = '0.918'
$c.Style = "Normal"
$ws.Range("E45").Value = '  +9.20%  '
$ws.Range("D46").Value = '1.731.53'
$ws.Range("E46").Value = '  +0.24%  '
$c = $ws.Range("D47")
$c.NumberFormat = "@"
$c.Value = '89.68'
$c.Style = "Normal"
$ws.Range("E47").Value = '  -0.32%  '
$ws.Range("E48").Value = '  -0.37%  '
$ws.Range("E49").Value = '  -1.46%  '
$ws.Range("E50").Value = '  +3.45%  '
$ws.Range("E51").Value = '  -1.43%  '
